$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename header strings: "<field>_old" -> "<field>_FV2404", "<field>_new" -> "<field>_FV2410" ---
$fields = @("Segmentname", "Segmentgruppe", "Segment", "Datenelement", "Segment ID", "Code", "Qualifier", "Beschreibung", "Bedingungsausdruck", "Bedingung")

for ($i = 0; $i -lt $fields.Length; $i++) {
    $colLetter = [string][char](65 + $i)   # A..J
    $ws.Range($colLetter + "1").Value = $fields[$i] + "_FV2404"
}
for ($i = 0; $i -lt $fields.Length; $i++) {
    $colLetter = [string][char](76 + $i)   # L..U
    $ws.Range($colLetter + "1").Value = $fields[$i] + "_FV2410"
}

# --- Freeze header row ---
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# --- Add a table over the full data range ---
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U60"), [System.Reflection.Missing]::Value, 1)
$tbl.Name = "Table1"

Write-Host "edit applied"
